$d = $word.ActiveDocument

# --- Edit 1: Paragraph 1 -- append a red "(This is a change ... )" note ---
$p1 = $d.Paragraphs(1)
$r = $p1.Range
$r.End = $r.End - 1          # exclude the paragraph mark
$r.InsertAfter("  ")          # two trailing spaces after the original sentence

$ins1 = $d.Range($r.End, $r.End)
$ins1.InsertAfter([char]0x0028 + "This is a change " + [char]0x2013 + " Ve")
$ins1.Font.Color = 192        # COLORREF 0x0000C0 => w:val="C00000"

$ins2 = $d.Range($ins1.End, $ins1.End)
$ins2.InsertAfter("rsion for branch alternate")
$ins2.Font.Color = 192

$ins3 = $d.Range($ins2.End, $ins2.End)
$ins3.InsertAfter(")")
$ins3.Font.Color = 192

# --- Edit 2: append a new, empty, shaded paragraph at the very end of the body ---
$endRange = $d.Range($d.Content.End, $d.Content.End)
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:pPr></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$endRange.InsertXML($xml)

# --- Housekeeping: drop styles that are no longer referenced anywhere in the
#     document (mirrors the style list pruning Word performs when it
#     re-saves a .docx). Deleted from the most-recently-added custom style
#     backwards to avoid destabilizing the Styles collection indices.       ---
$unusedStyles = @(
    "podcast-tools__subscribe-links",
    "generic-title",
    "subscribe-more-info",
    "subscribe",
    "audio-tool",
    "Heading 4 Char",
    "Heading 2 Char",
    "Hyperlink",
    "apple-converted-space",
    "Heading 4",
    "Heading 2"
)
foreach ($styleName in $unusedStyles) {
    $d.Styles($styleName).Delete()
}
